$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 180 - shifts existing rows 180..248 down to 181..249,
# matching the dimension change from A1:R248 to A1:R249.
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new data record.
$ws.Cells.Item(180, 1).Value2 = 6
$ws.Cells.Item(180, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(180, 3).Value2 = "Metropolitana"
$ws.Cells.Item(180, 4).Value2 = 44825
$ws.Cells.Item(180, 5).Value2 = 13
$ws.Cells.Item(180, 6).Value2 = 100112001
$ws.Cells.Item(180, 7).Value2 = "Berenjena"
$ws.Cells.Item(180, 8).Value2 = "Sin especificar"
$ws.Cells.Item(180, 9).Value2 = "Primera"
$ws.Cells.Item(180, 10).Value2 = 120
$ws.Cells.Item(180, 11).Value2 = 10000
$ws.Cells.Item(180, 12).Value2 = 11000
$ws.Cells.Item(180, 13).Value2 = 10583
$ws.Cells.Item(180, 14).Value2 = "`$/caja 40 unidades"
$ws.Cells.Item(180, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(180, 16).Value2 = 265
$ws.Cells.Item(180, 17).Value2 = 40
$ws.Cells.Item(180, 18).Value2 = "Hortaliza"

# Match the date-time number format used by the rest of column D.
$ws.Cells.Item(180, 4).NumberFormat = $ws.Cells.Item(181, 4).NumberFormat
